$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10.85093659812083
$ws.Range("D2").Value = 6.677977662816765
$ws.Range("E2").Value = 12.19200253495039
$ws.Range("F2").Value = 38.3689717592454
$ws.Range("G2").Value = 3.701182552676575
$ws.Range("J2").Value = 9.876820218374377
$ws.Range("L2").Value = 8.979474351254776
$ws.Range("M2").Value = 61.04346289703845
$ws.Range("O2").Value = 30.55881363894231
$ws.Range("C3").Value = 11.03237320566248
$ws.Range("D3").Value = 6.704475832012226
$ws.Range("E3").Value = 12.22349269785176
$ws.Range("F3").Value = 38.89944246664143
$ws.Range("G3").Value = 3.705739107834648
$ws.Range("J3").Value = 9.934950772748087
$ws.Range("L3").Value = 8.945956109344472
$ws.Range("M3").Value = 57.78014407461364
$ws.Range("O3").Value = 30.92669714220865
$ws.Range("C4").Value = 11.14787825913656
$ws.Range("D4").Value = 6.721779562582117
$ws.Range("E4").Value = 12.24749840947298
$ws.Range("F4").Value = 39.2428189464687
$ws.Range("G4").Value = 3.708655431829699
$ws.Range("J4").Value = 9.973786076077483
$ws.Range("L4").Value = 8.927282661367455
$ws.Range("M4").Value = 55.67108106093676
$ws.Range("O4").Value = 31.16589364671878
$ws.Range("C5").Value = 11.19598728441803
$ws.Range("D5").Value = 6.729091289936149
$ws.Range("E5").Value = 12.25843941746869
$ws.Range("F5").Value = 39.38715268302381
$ws.Range("G5").Value = 3.709873864391372
$ws.Range("J5").Value = 9.990393645488023
$ws.Range("L5").Value = 8.920155023320271
$ws.Range("M5").Value = 54.78554212256898
$ws.Range("O5").Value = 31.26668455027256
$ws.Range("C6").Value = 11.20403878801677
$ws.Range("D6").Value = 6.730321130532606
$ws.Range("E6").Value = 12.26032566906076
$ws.Range("F6").Value = 39.41138420818683
$ws.Range("G6").Value = 3.710078001975291
$ws.Range("J6").Value = 9.993198303835744
$ws.Range("L6").Value = 8.919000658831628
$ws.Range("M6").Value = 54.63693819973637
$ws.Range("O6").Value = 31.283620212254
$ws.Range("C7").Value = 11.14852285348616
$ws.Range("D7").Value = 6.721877116561968
$ws.Range("E7").Value = 12.24764129403718
$ws.Range("F7").Value = 39.24474769286181
$ws.Range("G7").Value = 3.708671742292358
$ws.Range("J7").Value = 9.974006896638327
$ws.Range("L7").Value = 8.927184581133936
$ws.Range("M7").Value = 55.6592433382386
$ws.Range("O7").Value = 31.16723956062055
$ws.Range("C8").Value = 10.91264957285821
$ws.Range("D8").Value = 6.686899896060862
$ws.Range("E8").Value = 12.20188128833297
$ws.Range("F8").Value = 38.54818858397164
$ws.Range("G8").Value = 3.702729156904813
$ws.Range("J8").Value = 9.896206117143443
$ws.Range("L8").Value = 8.967523018672404
$ws.Range("M8").Value = 59.94036112502465
$ws.Range("O8").Value = 30.68287687314393
$ws.Range("C9").Value = 10.48226083249066
$ws.Range("D9").Value = 6.626496044217394
$ws.Range("E9").Value = 12.1499313605829
$ws.Range("F9").Value = 37.32418599202199
$ws.Range("G9").Value = 3.692007756141586
$ws.Range("J9").Value = 9.768976445368319
$ws.Range("L9").Value = 9.061643151019661
$ws.Range("M9").Value = 67.4858687996865
$ws.Range("O9").Value = 29.8401896316674
$ws.Range("C10").Value = 10.18508723593381
$ws.Range("D10").Value = 6.58708728330239
$ws.Range("E10").Value = 12.13580782904986
$ws.Range("F10").Value = 36.51410599770928
$ws.Range("G10").Value = 3.684686038604067
$ws.Range("J10").Value = 9.691529834261377
$ws.Range("L10").Value = 9.139788009473426
$ws.Range("M10").Value = 72.49861629576989
$ws.Range("O10").Value = 29.28859118095484
$ws.Range("C11").Value = 10.05389865413492
$ws.Range("D11").Value = 6.57023458516864
$ws.Range("E11").Value = 12.1348226699101
$ws.Range("F11").Value = 36.16560255650532
$ws.Range("G11").Value = 3.681472890672149
$ws.Range("J11").Value = 9.659914088725476
$ws.Range("L11").Value = 9.177253611322687
$ws.Range("M11").Value = 74.66225675096003
$ws.Range("O11").Value = 29.05284556673216
$ws.Range("C12").Value = 10.004785115916
$ws.Range("D12").Value = 6.564007200709332
$ws.Range("E12").Value = 12.13524983604026
$ws.Range("F12").Value = 36.036572031167
$ws.Range("G12").Value = 3.680272822341397
$ws.Range("J12").Value = 9.648473939284145
$ws.Range("L12").Value = 9.191712946990114
$ws.Range("M12").Value = 75.46472367327593
$ws.Range("O12").Value = 28.96580620621261
$ws.Range("C13").Value = 10.01533764252862
$ws.Range("D13").Value = 6.565341515918844
$ws.Range("E13").Value = 12.13512196081215
$ws.Range("F13").Value = 36.06422923200173
$ws.Range("G13").Value = 3.680530540218135
$ws.Range("J13").Value = 9.650913916500963
$ws.Range("L13").Value = 9.188586837363751
$ws.Range("M13").Value = 75.29264861648062
$ws.Range("O13").Value = 28.98445157652805
$ws.Range("C14").Value = 10.04984679469015
$ws.Range("D14").Value = 6.569719161533179
$ws.Range("E14").Value = 12.13484167222908
$ws.Range("F14").Value = 36.15492783026011
$ws.Range("G14").Value = 3.681373827304423
$ws.Range("J14").Value = 9.658962165757069
$ws.Range("L14").Value = 9.178437761733798
$ws.Range("M14").Value = 74.72861433494394
$ws.Range("O14").Value = 29.04563974932888
$ws.Range("C15").Value = 10.07105788676726
$ws.Range("D15").Value = 6.57242069502562
$ws.Range("E15").Value = 12.13477473525004
$ws.Range("F15").Value = 36.2108682089673
$ws.Range("G15").Value = 3.681892530695141
$ws.Range("J15").Value = 9.663961611376621
$ws.Range("L15").Value = 9.172256459121865
$ws.Range("M15").Value = 74.38093013600401
$ws.Range("O15").Value = 29.08341152637281
$ws.Range("C16").Value = 10.19374029983226
$ws.Range("D16").Value = 6.588210255582047
$ws.Range("E16").Value = 12.13598320986268
$ws.Range("F16").Value = 36.53728937959742
$ws.Range("G16").Value = 3.684898371817855
$ws.Range("J16").Value = 9.69366977185406
$ws.Range("L16").Value = 9.137377803454275
$ws.Range("M16").Value = 72.35486180454296
$ws.Range("O16").Value = 29.30430732118314
$ws.Range("C17").Value = 10.27001871262121
$ws.Range("D17").Value = 6.598171737281564
$ws.Range("E17").Value = 12.13813059882677
$ws.Range("F17").Value = 36.74270558294314
$ws.Range("G17").Value = 3.686772307228887
$ws.Range("J17").Value = 9.712829171247733
$ws.Range("L17").Value = 9.116469266308162
$ws.Range("M17").Value = 71.08198128347144
$ws.Range("O17").Value = 29.44374174842028
$ws.Range("C18").Value = 10.31426889430066
$ws.Range("D18").Value = 6.604002466290604
$ws.Range("E18").Value = 12.1398766503611
$ws.Range("F18").Value = 36.86273450506147
$ws.Range("G18").Value = 3.687861222147885
$ws.Range("J18").Value = 9.72418865538549
$ws.Range("L18").Value = 9.104623910277777
$ws.Range("M18").Value = 70.33887322576025
$ws.Range("O18").Value = 29.52536631717533
$ws.Range("C19").Value = 10.32931627167602
$ws.Range("D19").Value = 6.605994028438119
$ws.Range("E19").Value = 12.14055503856893
$ws.Range("F19").Value = 36.90369510592012
$ws.Range("G19").Value = 3.688231819109154
$ws.Range("J19").Value = 9.728092729834447
$ws.Range("L19").Value = 9.100644420800508
$ws.Range("M19").Value = 70.08538664006385
$ws.Range("O19").Value = 29.55324638099002
$ws.Range("C20").Value = 10.26185981394673
$ws.Range("D20").Value = 6.597100853306756
$ws.Range("E20").Value = 12.13784899605838
$ws.Range("F20").Value = 36.72064377122845
$ws.Range("G20").Value = 3.686571678774505
$ws.Range("J20").Value = 9.710754396620571
$ws.Range("L20").Value = 9.118676342738476
$ws.Range("M20").Value = 71.21861874814857
$ws.Range("O20").Value = 29.42875082559502
$ws.Range("C21").Value = 10.03969537542136
$ws.Range("D21").Value = 6.568429152450214
$ws.Range("E21").Value = 12.1349021381421
$ws.Range("F21").Value = 36.12820712619158
$ws.Range("G21").Value = 3.681125682517354
$ws.Range("J21").Value = 9.656583659972872
$ws.Range("L21").Value = 9.181411441737023
$ws.Range("M21").Value = 74.89474279860987
$ws.Range("O21").Value = 29.02760630058198
$ws.Range("C22").Value = 9.897784832582039
$ws.Range("D22").Value = 6.550590260326573
$ws.Range("E22").Value = 12.13764890002095
$ws.Range("F22").Value = 35.75818643900286
$ws.Range("G22").Value = 3.677663532020582
$ws.Range("J22").Value = 9.624287225158067
$ws.Range("L22").Value = 9.223995630543861
$ws.Range("M22").Value = 77.19903728923386
$ws.Range("O22").Value = 28.77847277497163
$ws.Range("C23").Value = 9.973227836748574
$ws.Range("D23").Value = 6.560028928944016
$ws.Range("E23").Value = 12.13574938929181
$ws.Range("F23").Value = 35.95407993719009
$ws.Range("G23").Value = 3.679502534775061
$ws.Range("J23").Value = 9.64123593585335
$ws.Range("L23").Value = 9.201124054754537
$ws.Range("M23").Value = 75.97819953868591
$ws.Range("O23").Value = 28.91022977431058
$ws.Range("C24").Value = 10.26554721562984
$ws.Range("D24").Value = 6.597584676811518
$ws.Range("E24").Value = 12.13797471685872
$ws.Range("F24").Value = 36.73061190200202
$ws.Range("G24").Value = 3.686662346861592
$ws.Range("J24").Value = 9.71169133013565
$ws.Range("L24").Value = 9.117677977365812
$ws.Range("M24").Value = 71.15688017870772
$ws.Range("O24").Value = 29.4355236721529
$ws.Range("C25").Value = 10.59530757797503
$ws.Range("D25").Value = 6.641962726132007
$ws.Range("E25").Value = 12.15983628627385
$ws.Range("F25").Value = 37.63985457578025
$ws.Range("G25").Value = 3.69480971287985
$ws.Range("J25").Value = 9.800624709892075
$ws.Range("L25").Value = 9.034589161035884
$ws.Range("M25").Value = 65.53716511521169
$ws.Range("O25").Value = 30.05647015105012
